$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently wraps the
#    image paragraph. It gets re-created further up the document
#    (inside the new page-break paragraph) in step 3 below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Extend the "method that took the most time" paragraph with the
#    additional sentences describing the Ruby multithreading issue
#    and the hash-caching optimization.
# ------------------------------------------------------------------
$oldText = "The method that took the most time was related to verifying the hash. Our initial implementation tried to remedy a bit of the pain associated with this method by using hashes instead of lists and to try not to repeat string manipulations. We then attempted to use multithreading to run somewhat costly processes simultaneously to cut down on runtime."
$newText = "The method that took the most time was related to verifying the hash. Our initial implementation tried to remedy a bit of the pain associated with this method by using hashes instead of lists and to try not to repeat string manipulations. We then attempted to use multithreading to run somewhat costly processes simultaneously to cut down on runtime, but this ended up being tricky to implement in Ruby. We then switched to storing the calculated hash value of each character so that we didn't need to re-calculate hash values that we have already calculated. This change presented the greatest time improvement."

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# ------------------------------------------------------------------
# 3) Collapse the four empty paragraphs that followed into a single
#    paragraph containing a page break, carrying the "_GoBack"
#    bookmark (both start and end) inside it.
# ------------------------------------------------------------------

# Find the paragraph that holds the just-edited text (still unique),
# identified by index; the paragraph right after it is the first of
# the four empty paragraphs that must collapse into the new break
# paragraph.
$editedIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.IndexOf("This change presented the greatest time improvement.") -ge 0) {
        $editedIndex = $i
        break
    }
}

$afterEditIndex = $editedIndex + 1
$afterEditPara = $d.Paragraphs.Item($afterEditIndex)
$insertPoint = $d.Range($afterEditPara.Range.Start, $afterEditPara.Range.Start)

$breakXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:br w:type="page"/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($breakXml)

# Remove the now-redundant paragraphs: the placeholder empty
# paragraph minted by InsertXML plus the four original empty
# paragraphs, leaving only the single page-break paragraph behind.
$firstToRemoveIndex = $afterEditIndex + 1
$lastToRemoveIndex = $firstToRemoveIndex + 4
$firstToRemove = $d.Paragraphs.Item($firstToRemoveIndex)
$lastToRemove = $d.Paragraphs.Item($lastToRemoveIndex)
$removeRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$removeRange.Delete()
